$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at position 5. This shifts the existing rows 5-135
#    down by one (old row 5 -> new row 6, ..., old row 13 -> new row 14).
#    The new row 5 is populated below with the data that (per the diff)
#    belongs to record "A 38880-2020" with refreshed survey figures; the
#    duplicate left behind at (new) row 14 is removed afterwards.
# ---------------------------------------------------------------------------
$ws.Rows(5).Insert()

$ws.Cells.Item(5, 1).Value = "A 38880-2020"
$ws.Cells.Item(5, 2).Value2 = 44062
$ws.Cells.Item(5, 3).Value2 = 45192
$ws.Cells.Item(5, 4).Value = "SÖDERMANLANDS LÄN"
$ws.Cells.Item(5, 5).Value = "TROSA"
$ws.Cells.Item(5, 7).Value2 = 8.9
$ws.Cells.Item(5, 8).Value2 = 4
$ws.Cells.Item(5, 9).Value2 = 0
$ws.Cells.Item(5, 10).Value2 = 0
$ws.Cells.Item(5, 11).Value2 = 1
$ws.Cells.Item(5, 12).Value2 = 0
$ws.Cells.Item(5, 13).Value2 = 0
$ws.Cells.Item(5, 14).Value2 = 0
$ws.Cells.Item(5, 15).Value2 = 1
$ws.Cells.Item(5, 16).Value2 = 1
$ws.Cells.Item(5, 17).Value2 = 4
$ws.Cells.Item(5, 18).Value = "Kricka`nÅkergroda`nMindre vattensalamander`nVanlig padda"

$ws.Cells.Item(5, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_TROSA/artfynd/A 38880-2020.xlsx", "A 38880-2020")'
$ws.Cells.Item(5, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_TROSA/kartor/A 38880-2020.png", "A 38880-2020")'
$ws.Cells.Item(5, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_TROSA/klagomål/A 38880-2020.docx", "A 38880-2020")'
$ws.Cells.Item(5, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_TROSA/klagomålsmail/A 38880-2020.docx", "A 38880-2020")'
$ws.Cells.Item(5, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_TROSA/tillsyn/A 38880-2020.docx", "A 38880-2020")'
$ws.Cells.Item(5, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_TROSA/tillsynsmail/A 38880-2020.docx", "A 38880-2020")'

$ws.Rows(5).RowHeight = 15

# ---------------------------------------------------------------------------
# 2) Remove the now-duplicated old row (original "A 38880-2020" record),
#    which the insert above pushed down to row 14.
# ---------------------------------------------------------------------------
$ws.Rows(14).Delete()

# ---------------------------------------------------------------------------
# 3) The "Förändrad" column (C) moves from 45190 to 45192 for every data
#    row (2 through 135).
# ---------------------------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
